$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", formatted like the existing header row (copy
# formatting from G1 "sum" so it keeps the same bold/centered/bordered style).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data cell H2 = 1 (plain number, no special formatting, like neighbors).
$ws.Range("H2").Value = 1
